$d = $word.ActiveDocument

# "In a nutshell, you create a project and write goals to it." ->
# "In a nutshell, you create a project and write down goals."
$d.Content.Find.Execute("write goals to it.", $false, $false, $false, $false, $false, $true, 1, $false, "write down goals.", 2)

# "You write notes about achieving the goals" ->
# "You write notes about ways of achieving the goals"
$d.Content.Find.Execute("notes about achieving", $false, $false, $false, $false, $false, $true, 1, $false, "notes about ways of achieving", 2)

# "You can open the associated e-book file in different pages," ->
# "You can open the associated e-book file at different bookmarked pages,"
$d.Content.Find.Execute("e-book file in different pages", $false, $false, $false, $false, $false, $true, 1, $false, "e-book file at different bookmarked pages", 2)

# "add tags to them and open them in a viewer." ->
# "add tags to them, and open them in a viewer."
$d.Content.Find.Execute("add tags to them and open", $false, $false, $false, $false, $false, $true, 1, $false, "add tags to them, and open", 2)
